# Apply the "Add files via upload" refresh: the Power-Query-backed table on
# Sheet1 picked up 26 additional rows of daily NAV data (rows 892-917),
# extending the query table, the workbook-level defined name, and the
# sheet dimension/view accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows: row number, Date serial, Patrimonio Neto (column G)
$newData = @(
    @(892, 44980, 11625.162),
    @(893, 44981, 11543.844999999999),
    @(894, 44984, 11601.199000000001),
    @(895, 44985, 11602.903),
    @(896, 44986, 11510.482),
    @(897, 44987, 11495.031999999999),
    @(898, 44988, 11566.897000000001),
    @(899, 44991, 11664.322),
    @(900, 44992, 11684.592000000001),
    @(901, 44993, 11766.861999999999),
    @(902, 44994, 11844.987999999999),
    @(903, 44995, 11786.33),
    @(904, 44998, 11784.51),
    @(905, 44999, 11820.071),
    @(906, 45000, 11724.504000000001),
    @(907, 45001, 11899.919),
    @(908, 45002, 11954.635),
    @(909, 45005, 11936.447),
    @(910, 45006, 11920.163),
    @(911, 45007, 11928.761),
    @(912, 45008, 11764.855),
    @(913, 45012, 11845.514999999999),
    @(914, 45013, 11957.04),
    @(915, 45014, 12145.429),
    @(916, 45015, 12411.696),
    @(917, 45016, 12470.038)
)

$lastRow = 917

# 1. Grow the query table (ListObject) so the new rows become part of it.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H$lastRow"))

# 2. Fill in every new row: Date, Fondo, asset_class, Horizonte, Patrimonio
#    Neto and the calculated "Valor CP" column (formula copied down like the
#    rest of the table).
foreach ($row in $newData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]                       # A: Date
    $ws.Cells.Item($r, 2).Value2 = "Adcap Balanceado - Clase B"  # B: Fondo Común de Inversión
    $ws.Cells.Item($r, 3).Value2 = "Retorno Total"                # C: asset_class
    $ws.Cells.Item($r, 5).Value2 = "FLEX"                         # E: Horizonte
    $ws.Cells.Item($r, 7).Value2 = $row[2]                        # G: Patrimonio Neto
    $ws.Cells.Item($r, 8).Formula = "=+Sheet1[[#This Row],[Valor mil Cuotapartes]]/1000"  # H: Valor CP
}

# 3. Match formatting: column A keeps the short-date style used by the rest
#    of the table, reuse the existing style (don't invent a new one).
$ws.Range("A891").Copy() | Out-Null
$ws.Range("A892:A$lastRow").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# The "Valor CP" column had picked up a spurious extra style (s="2") on the
# tail of the table; the refreshed table no longer carries it, so reset the
# whole calculated column back to the plain/default style.
$ws.Range("G891").Copy() | Out-Null
$ws.Range("H878:H$lastRow").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 4. Update the hidden Power Query defined name to cover the new data range.
$wb.Names.Item("DatosExternos_1").RefersTo = "=Sheet1!`$A`$1:`$G`$$lastRow"

# 5. Reflect the new extent/scroll position like Excel does after a refresh.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 892
$win.ScrollColumn = 1
$ws.Range("G922").Select()
